$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the leading Title / Author / Date paragraphs (first 3 paragraphs
#    of the document body: "Luke Richardson Resume" / "Luke Richardson" /
#    "May 2025").
# ---------------------------------------------------------------------------
$titleStart = $d.Paragraphs.Item(1).Range.Start
$dateEnd    = $d.Paragraphs.Item(3).Range.End
$d.Range($titleStart, $dateEnd).Delete()

# ---------------------------------------------------------------------------
# 2. Replace the "Experience" table (Where / When / What / Why) with a
#    single FirstParagraph-styled paragraph whose rows are separated by
#    manual line breaks, mirroring a markdown table that got pasted in as
#    plain text + hyperlinks.
# ---------------------------------------------------------------------------

# Locate the "Experience" heading paragraph that immediately precedes the
# table.
$expIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq "Experience") {
        $expIndex = $i
    }
}

$rows = @(
    @{ Where = "Camden";         Url = "https://www.camden.gov.uk/";           When = "2024 - 2025"; What = "Lead Network Engineer";    Why = "Public Service" },
    @{ Where = "Lloret";         Url = "https://www.lloret.co.uk/";            When = "2023 - 2024"; What = "Network Architect";        Why = "Rediscover my Roots" },
    @{ Where = "WeWork";         Url = "https://www.wework.com/";              When = "2019 - 2023"; What = "Global Network Architect"; Why = "Build Complex Systems" },
    @{ Where = "Redstone";       Url = "https://www.onnecgroup.com/";          When = "2017 - 2018"; What = "Network Engineer & TPM";    Why = "Prove Myself" },
    @{ Where = "Dimension Data"; Url = "https://www.dimensiondata.com/en-gb/"; When = "2012 - 2017"; What = "PM to Network Engineer";    Why = "Learn the Ropes" }
)

# ":– | :– | :– | :– |"  (the dash is U+2013 EN DASH, as produced by Word's
# autocorrect turning the markdown table rule "--" into an en dash).
$dash = [string][char]0x2013
$sep = ":" + $dash + " | :" + $dash + " | :" + $dash + " | :" + $dash + " |"

# Create the first (separator) paragraph right after "Experience".
$anchor = $d.Paragraphs.Item($expIndex)
$anchor.Range.InsertParagraphAfter()
$firstParaIndex = $expIndex + 1
$firstPara = $d.Paragraphs.Item($firstParaIndex)
$firstPara.Range.Style = "FirstParagraph"
$insertPoint = $d.Range($firstPara.Range.Start, $firstPara.Range.Start)
$insertPoint.InsertAfter($sep)

$prevParaIndex = $firstParaIndex

foreach ($row in $rows) {
    # Start a fresh paragraph for this row. Hyperlinks.Add always anchors to
    # the start of the paragraph containing its target Range, so the link
    # has to be the very first thing typed into a brand-new paragraph.
    $prevPara = $d.Paragraphs.Item($prevParaIndex)
    $endOfPrev = $prevPara.Range
    $endOfPrev.Collapse(0)
    $endOfPrev.InsertParagraphAfter()

    $rowParaIndex = $prevParaIndex + 1
    $rowPara = $d.Paragraphs.Item($rowParaIndex)
    $rowPara.Range.Style = "FirstParagraph"

    $hLoc = $d.Range($rowPara.Range.Start, $rowPara.Range.Start)
    [void]$d.Hyperlinks.Add($hLoc, $row.Url, "", "", $row.Where)

    $rowPara = $d.Paragraphs.Item($rowParaIndex)
    $tail = $rowPara.Range
    $tail.Collapse(0)
    $tail.InsertAfter(" | " + $row.When + " | " + $row.What + " | " + $row.Why)

    # Merge this row paragraph into the previous one: delete the paragraph
    # mark that separates them and drop a manual line break in its place.
    $prevPara = $d.Paragraphs.Item($prevParaIndex)
    $markPos = $prevPara.Range.End - 1
    $d.Range($markPos, $markPos + 1).Delete()
    $d.Range($markPos, $markPos).InsertAfter([string][char]11)
}

# Finally, delete the now-redundant table.
$d.Tables.Item(1).Delete()
